$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is forced to Text format before the write so that
# numeric-looking strings (e.g. "4.40", "18.90") keep their exact
# literal text instead of being auto-coerced to a number (which would
# drop significant trailing zeros). The style is reset back to the
# default "Normal" afterwards so no stray cell format/style is left
# behind (matching the original workbook, where these cells carry no
# explicit style).
$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '58.997.28'
$cell.Style = "Normal"
$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  -2.40%  '
$cell.Style = "Normal"
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.661.48'
$cell.Style = "Normal"
$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  -0.92%  '
$cell.Style = "Normal"
$cell = $ws.Range('E4')
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell.Style = "Normal"
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '525.27'
$cell.Style = "Normal"
$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  +0.51%  '
$cell.Style = "Normal"
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '144.32'
$cell.Style = "Normal"
$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  -1.31%  '
$cell.Style = "Normal"
$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  +0.27%  '
$cell.Style = "Normal"
$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  -1.01%  '
$cell.Style = "Normal"
$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +7.92%  '
$cell.Style = "Normal"
$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  -2.30%  '
$cell.Style = "Normal"
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.335'
$cell.Style = "Normal"
$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  -1.93%  '
$cell.Style = "Normal"
$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  +1.39%  '
$cell.Style = "Normal"
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '3.130.05'
$cell.Style = "Normal"
$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  -0.87%  '
$cell.Style = "Normal"
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '58.995.05'
$cell.Style = "Normal"
$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  -2.41%  '
$cell.Style = "Normal"
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '21.06'
$cell.Style = "Normal"
$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  -0.96%  '
$cell.Style = "Normal"
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '2.670.65'
$cell.Style = "Normal"
$cell = $ws.Range('E16')
$cell.NumberFormat = "@"
$cell.Value = '  -3.15%  '
$cell.Style = "Normal"
$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  -1.80%  '
$cell.Style = "Normal"
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '338.58'
$cell.Style = "Normal"
$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  -3.54%  '
$cell.Style = "Normal"
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '4.40'
$cell.Style = "Normal"
$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  -3.18%  '
$cell.Style = "Normal"
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '10.39'
$cell.Style = "Normal"
$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  -1.64%  '
$cell.Style = "Normal"
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '6.41'
$cell.Style = "Normal"
$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +1.28%  '
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.Style = "Normal"
$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  +2.36%  '
$cell.Style = "Normal"
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '0.418'
$cell.Style = "Normal"
$cell = $ws.Range('E24')
$cell.NumberFormat = "@"
$cell.Value = '  -0.81%  '
$cell.Style = "Normal"
$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  -1.69%  '
$cell.Style = "Normal"
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$cell = $ws.Range('E26')
$cell.NumberFormat = "@"
$cell.Value = '  +0.29%  '
$cell.Style = "Normal"
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0800'
$cell.Style = "Normal"
$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  -1.85%  '
$cell.Style = "Normal"
$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  -2.50%  '
$cell.Style = "Normal"
$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  -2.63%  '
$cell.Style = "Normal"
$cell = $ws.Range('E30')
$cell.NumberFormat = "@"
$cell.Value = '  +0.10%  '
$cell.Style = "Normal"
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '1.60'
$cell.Style = "Normal"
$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  -0.20%  '
$cell.Style = "Normal"
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '18.85'
$cell.Style = "Normal"
$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  -1.08%  '
$cell.Style = "Normal"
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '150.54'
$cell.Style = "Normal"
$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  +1.47%  '
$cell.Style = "Normal"
$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  -3.97%  '
$cell.Style = "Normal"
$cell = $ws.Range('E35')
$cell.NumberFormat = "@"
$cell.Value = '  -3.84%  '
$cell.Style = "Normal"
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.892'
$cell.Style = "Normal"
$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  -6.04%  '
$cell.Style = "Normal"
$cell = $ws.Range('E37')
$cell.NumberFormat = "@"
$cell.Value = '  -0.59%  '
$cell.Style = "Normal"
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '36.89'
$cell.Style = "Normal"
$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  +0.20%  '
$cell.Style = "Normal"
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.46'
$cell.Style = "Normal"
$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  -6.09%  '
$cell.Style = "Normal"
$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  -3.07%  '
$cell.Style = "Normal"
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.616'
$cell.Style = "Normal"
$cell = $ws.Range('E41')
$cell.NumberFormat = "@"
$cell.Value = '  +0.48%  '
$cell.Style = "Normal"
$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  +0.36%  '
$cell.Style = "Normal"
$cell = $ws.Range('B43')
$cell.NumberFormat = "@"
$cell.Value = 'Bittensor'
$cell.Style = "Normal"
$cell = $ws.Range('C43')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell.Style = "Normal"
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '275.57'
$cell.Style = "Normal"
$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  -2.47%  '
$cell.Style = "Normal"
$cell = $ws.Range('B44')
$cell.NumberFormat = "@"
$cell.Value = 'EnergySwap'
$cell.Style = "Normal"
$cell = $ws.Range('C44')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.Style = "Normal"
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '19.88'
$cell.Style = "Normal"
$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  -0.35%  '
$cell.Style = "Normal"
$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  -1.73%  '
$cell.Style = "Normal"
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '10.67'
$cell.Style = "Normal"
$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  +2.05%  '
$cell.Style = "Normal"
$cell = $ws.Range('B47')
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell.Style = "Normal"
$cell = $ws.Range('C47')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell.Style = "Normal"
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '2.048.58'
$cell.Style = "Normal"
$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  -3.71%  '
$cell.Style = "Normal"
$cell = $ws.Range('B48')
$cell.NumberFormat = "@"
$cell.Value = 'Hedera'
$cell.Style = "Normal"
$cell = $ws.Range('C48')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell.Style = "Normal"
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.0530'
$cell.Style = "Normal"
$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  -1.55%  '
$cell.Style = "Normal"
$cell = $ws.Range('B49')
$cell.NumberFormat = "@"
$cell.Value = 'RenderToken'
$cell.Style = "Normal"
$cell = $ws.Range('C49')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.Style = "Normal"
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '4.71'
$cell.Style = "Normal"
$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  -3.28%  '
$cell.Style = "Normal"
$cell = $ws.Range('B50')
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell.Style = "Normal"
$cell = $ws.Range('C50')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.Style = "Normal"
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.0229'
$cell.Style = "Normal"
$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  -2.66%  '
$cell.Style = "Normal"
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '18.90'
$cell.Style = "Normal"
$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  -1.45%  '
$cell.Style = "Normal"
